# This script applies updated NATMI ligand-receptor interaction metrics
# for the Lpl-Lrp1 pair (Sheet1), following Dr Hou's advice.
# Ligand-expressing cells (E) and Receptor-expressing cells (K) counts
# increase from 1 to 3 for every data row, with all dependent average/
# total expression and specificity metrics (columns G-J, M-T) recomputed
# accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 46.48074800000001
$ws.Cells.Item(2, 8).Value = 139.442244
$ws.Cells.Item(2, 9).Value = 0.1473944418036112
$ws.Cells.Item(2, 10).Value = 0.1473944418036112
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 23.63579766666667
$ws.Cells.Item(2, 14).Value = 70.907393
$ws.Cells.Item(2, 15).Value = 0.06827844587621175
$ws.Cells.Item(2, 16).Value = 0.06827844587621175
$ws.Cells.Item(2, 17).Value = 1098.609555123321
$ws.Cells.Item(2, 18).Value = 9887.485996109894
$ws.Cells.Item(2, 19).Value = 0.01006386341714231
$ws.Cells.Item(2, 20).Value = 0.01006386341714231
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 46.48074800000001
$ws.Cells.Item(3, 8).Value = 139.442244
$ws.Cells.Item(3, 9).Value = 0.1473944418036112
$ws.Cells.Item(3, 10).Value = 0.1473944418036112
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 181.2883913333334
$ws.Cells.Item(3, 14).Value = 543.865174
$ws.Cells.Item(3, 15).Value = 0.5237009467675041
$ws.Cells.Item(3, 16).Value = 0.523700946767504
$ws.Cells.Item(3, 17).Value = 8426.420032890053
$ws.Cells.Item(3, 18).Value = 75837.78029601047
$ws.Cells.Item(3, 19).Value = 0.07719060872081895
$ws.Cells.Item(3, 20).Value = 0.07719060872081894
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 46.48074800000001
$ws.Cells.Item(4, 8).Value = 139.442244
$ws.Cells.Item(4, 9).Value = 0.1473944418036112
$ws.Cells.Item(4, 10).Value = 0.1473944418036112
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 111.1005463333333
$ws.Cells.Item(4, 14).Value = 333.301639
$ws.Cells.Item(4, 15).Value = 0.3209442197221123
$ws.Cells.Item(4, 16).Value = 0.3209442197221123
$ws.Cells.Item(4, 17).Value = 5164.036496781991
$ws.Cells.Item(4, 18).Value = 46476.32847103792
$ws.Cells.Item(4, 19).Value = 0.04730539411603627
$ws.Cells.Item(4, 20).Value = 0.04730539411603627
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 46.48074800000001
$ws.Cells.Item(5, 8).Value = 139.442244
$ws.Cells.Item(5, 9).Value = 0.1473944418036112
$ws.Cells.Item(5, 10).Value = 0.1473944418036112
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 30.14303933333333
$ws.Cells.Item(5, 14).Value = 90.42911799999999
$ws.Cells.Item(5, 15).Value = 0.08707638763417187
$ws.Cells.Item(5, 16).Value = 0.08707638763417187
$ws.Cells.Item(5, 17).Value = 1401.071015206755
$ws.Cells.Item(5, 18).Value = 12609.63913686079
$ws.Cells.Item(5, 19).Value = 0.01283457554961363
$ws.Cells.Item(5, 20).Value = 0.01283457554961363
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 192.830597
$ws.Cells.Item(6, 8).Value = 578.4917909999999
$ws.Cells.Item(6, 9).Value = 0.6114823756165045
$ws.Cells.Item(6, 10).Value = 0.6114823756165044
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 23.63579766666667
$ws.Cells.Item(6, 14).Value = 70.907393
$ws.Cells.Item(6, 15).Value = 0.06827844587621175
$ws.Cells.Item(6, 16).Value = 0.06827844587621175
$ws.Cells.Item(6, 17).Value = 4557.704974634539
$ws.Cells.Item(6, 18).Value = 41019.34477171086
$ws.Cells.Item(6, 19).Value = 0.04175106628778889
$ws.Cells.Item(6, 20).Value = 0.04175106628778888
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 192.830597
$ws.Cells.Item(7, 8).Value = 578.4917909999999
$ws.Cells.Item(7, 9).Value = 0.6114823756165045
$ws.Cells.Item(7, 10).Value = 0.6114823756165044
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 181.2883913333334
$ws.Cells.Item(7, 14).Value = 543.865174
$ws.Cells.Item(7, 15).Value = 0.5237009467675041
$ws.Cells.Item(7, 16).Value = 0.523700946767504
$ws.Cells.Item(7, 17).Value = 34957.94872997629
$ws.Cells.Item(7, 18).Value = 314621.5385697866
$ws.Cells.Item(7, 19).Value = 0.320233899042006
$ws.Cells.Item(7, 20).Value = 0.3202338990420059
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 192.830597
$ws.Cells.Item(8, 8).Value = 578.4917909999999
$ws.Cells.Item(8, 9).Value = 0.6114823756165045
$ws.Cells.Item(8, 10).Value = 0.6114823756165044
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 111.1005463333333
$ws.Cells.Item(8, 14).Value = 333.301639
$ws.Cells.Item(8, 15).Value = 0.3209442197221123
$ws.Cells.Item(8, 16).Value = 0.3209442197221123
$ws.Cells.Item(8, 17).Value = 21423.58467648282
$ws.Cells.Item(8, 18).Value = 192812.2620883454
$ws.Cells.Item(8, 19).Value = 0.1962517339160626
$ws.Cells.Item(8, 20).Value = 0.1962517339160626
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 192.830597
$ws.Cells.Item(9, 8).Value = 578.4917909999999
$ws.Cells.Item(9, 9).Value = 0.6114823756165045
$ws.Cells.Item(9, 10).Value = 0.6114823756165044
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 30.14303933333333
$ws.Cells.Item(9, 14).Value = 90.42911799999999
$ws.Cells.Item(9, 15).Value = 0.08707638763417187
$ws.Cells.Item(9, 16).Value = 0.08707638763417187
$ws.Cells.Item(9, 17).Value = 5812.500270041148
$ws.Cells.Item(9, 18).Value = 52312.50243037032
$ws.Cells.Item(9, 19).Value = 0.05324567637064703
$ws.Cells.Item(9, 20).Value = 0.05324567637064703
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 60.030993
$ws.Cells.Item(10, 8).Value = 180.092979
$ws.Cells.Item(10, 9).Value = 0.1903634318482028
$ws.Cells.Item(10, 10).Value = 0.1903634318482028
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 23.63579766666667
$ws.Cells.Item(10, 14).Value = 70.907393
$ws.Cells.Item(10, 15).Value = 0.06827844587621175
$ws.Cells.Item(10, 16).Value = 0.06827844587621175
$ws.Cells.Item(10, 17).Value = 1418.880404277083
$ws.Cells.Item(10, 18).Value = 12769.92363849375
$ws.Cells.Item(10, 19).Value = 0.01299771927825744
$ws.Cells.Item(10, 20).Value = 0.01299771927825744
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 60.030993
$ws.Cells.Item(11, 8).Value = 180.092979
$ws.Cells.Item(11, 9).Value = 0.1903634318482028
$ws.Cells.Item(11, 10).Value = 0.1903634318482028
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 181.2883913333334
$ws.Cells.Item(11, 14).Value = 543.865174
$ws.Cells.Item(11, 15).Value = 0.5237009467675041
$ws.Cells.Item(11, 16).Value = 0.523700946767504
$ws.Cells.Item(11, 17).Value = 10882.9221511126
$ws.Cells.Item(11, 18).Value = 97946.29936001336
$ws.Cells.Item(11, 19).Value = 0.09969350948881506
$ws.Cells.Item(11, 20).Value = 0.09969350948881503
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 60.030993
$ws.Cells.Item(12, 8).Value = 180.092979
$ws.Cells.Item(12, 9).Value = 0.1903634318482028
$ws.Cells.Item(12, 10).Value = 0.1903634318482028
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 111.1005463333333
$ws.Cells.Item(12, 14).Value = 333.301639
$ws.Cells.Item(12, 15).Value = 0.3209442197221123
$ws.Cells.Item(12, 16).Value = 0.3209442197221123
$ws.Cells.Item(12, 17).Value = 6669.476119232509
$ws.Cells.Item(12, 18).Value = 60025.28507309258
$ws.Cells.Item(12, 19).Value = 0.06109604309814495
$ws.Cells.Item(12, 20).Value = 0.06109604309814495
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 60.030993
$ws.Cells.Item(13, 8).Value = 180.092979
$ws.Cells.Item(13, 9).Value = 0.1903634318482028
$ws.Cells.Item(13, 10).Value = 0.1903634318482028
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 30.14303933333333
$ws.Cells.Item(13, 14).Value = 90.42911799999999
$ws.Cells.Item(13, 15).Value = 0.08707638763417187
$ws.Cells.Item(13, 16).Value = 0.08707638763417187
$ws.Cells.Item(13, 17).Value = 1809.516583218058
$ws.Cells.Item(13, 18).Value = 16285.64924896252
$ws.Cells.Item(13, 19).Value = 0.01657615998298537
$ws.Cells.Item(13, 20).Value = 0.01657615998298537
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 16.00705666666667
$ws.Cells.Item(14, 8).Value = 48.02117
$ws.Cells.Item(14, 9).Value = 0.05075975073168155
$ws.Cells.Item(14, 10).Value = 0.05075975073168155
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 23.63579766666667
$ws.Cells.Item(14, 14).Value = 70.907393
$ws.Cells.Item(14, 15).Value = 0.06827844587621175
$ws.Cells.Item(14, 16).Value = 0.06827844587621175
$ws.Cells.Item(14, 17).Value = 378.3395526122011
$ws.Cells.Item(14, 18).Value = 3405.05597350981
$ws.Cells.Item(14, 19).Value = 0.003465796893023119
$ws.Cells.Item(14, 20).Value = 0.003465796893023118
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 16.00705666666667
$ws.Cells.Item(15, 8).Value = 48.02117
$ws.Cells.Item(15, 9).Value = 0.05075975073168155
$ws.Cells.Item(15, 10).Value = 0.05075975073168155
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 181.2883913333334
$ws.Cells.Item(15, 14).Value = 543.865174
$ws.Cells.Item(15, 15).Value = 0.5237009467675041
$ws.Cells.Item(15, 16).Value = 0.523700946767504
$ws.Cells.Item(15, 17).Value = 2901.893553081509
$ws.Cells.Item(15, 18).Value = 26117.04197773358
$ws.Cells.Item(15, 19).Value = 0.02658292951586414
$ws.Cells.Item(15, 20).Value = 0.02658292951586413
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 16.00705666666667
$ws.Cells.Item(16, 8).Value = 48.02117
$ws.Cells.Item(16, 9).Value = 0.05075975073168155
$ws.Cells.Item(16, 10).Value = 0.05075975073168155
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 111.1005463333333
$ws.Cells.Item(16, 14).Value = 333.301639
$ws.Cells.Item(16, 15).Value = 0.3209442197221123
$ws.Cells.Item(16, 16).Value = 0.3209442197221123
$ws.Cells.Item(16, 17).Value = 1778.392740855292
$ws.Cells.Item(16, 18).Value = 16005.53466769763
$ws.Cells.Item(16, 19).Value = 0.01629104859186846
$ws.Cells.Item(16, 20).Value = 0.01629104859186845
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 16.00705666666667
$ws.Cells.Item(17, 8).Value = 48.02117
$ws.Cells.Item(17, 9).Value = 0.05075975073168155
$ws.Cells.Item(17, 10).Value = 0.05075975073168155
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 30.14303933333333
$ws.Cells.Item(17, 14).Value = 90.42911799999999
$ws.Cells.Item(17, 15).Value = 0.08707638763417187
$ws.Cells.Item(17, 16).Value = 0.08707638763417187
$ws.Cells.Item(17, 17).Value = 1809.516583218058
$ws.Cells.Item(17, 18).Value = 16285.64924896252
$ws.Cells.Item(17, 19).Value = 0.004419975730925842
$ws.Cells.Item(17, 20).Value = 0.004419975730925842

Write-Output "Updated 224 cell values on Sheet1."
